$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values like "1.001" or "0.07210" stored as literal
# text (note some are multi-dot "thousands" strings like "26.065.54"). A plain
# .Value assignment lets Excel auto-parse number-shaped strings into real numbers
# (e.g. "1.001" -> 1, "0.07203" -> 7.203E-02), which corrupts the intended text.
# Force the whole column to Text format first so every literal is stored verbatim,
# then restore the default "Normal" style afterwards so no cell keeps a leftover
# explicit number format / style index.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.054.30"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "1.747.05"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "233.38"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.5263"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "40.13"
$ws.Range("E9").Value = "  +4.55%  "
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").Value = "1.755.30"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "0.07203"
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "0.6391"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "78.25"
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "25.985.96"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "11.55"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "0.000006689"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").Value = "1.976.27"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").Value = "4.321"
$ws.Range("E23").Value = "  +6.86%  "
$ws.Range("D24").Value = "8.810"
$ws.Range("E24").Value = "  +4.21%  "
$ws.Range("D25").Value = "5.191"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").Value = "139.53"
$ws.Range("E26").Value = "  +2.15%  "
$ws.Range("D27").Value = "1.522"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Value = "1.803"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").Value = "104.25"
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D31").Value = "0.08286"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").Value = "3.765"
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("D33").Value = "3.649"
$ws.Range("E33").Value = "  +8.13%  "
$ws.Range("D34").Value = "0.04518"
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").Value = "2.635"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").Value = "0.9966"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").Value = "0.6291"
$ws.Range("E37").Value = "  +5.43%  "
$ws.Range("D38").Value = "2.706"
$ws.Range("D39").Value = "0.01589"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").Value = "1.922"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").Value = "97.55"
$ws.Range("E42").Value = "  -4.25%  "
$ws.Range("D43").Value = "0.3883"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("D44").Value = "0.7328"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("D45").Value = "5.033"
$ws.Range("E45").Value = "  +3.34%  "
$ws.Range("D46").Value = "0.1140"
$ws.Range("E46").Value = "  +3.73%  "
$ws.Range("D47").Value = "0.05342"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").Value = "6.290"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "53.86"
$ws.Range("E49").Value = "  +3.85%  "
$ws.Range("D50").Value = "30.48"
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("D51").Value = "7.668"
$ws.Range("E51").Value = "  +3.47%  "

$priceRange.Style = "Normal"
